$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" column (N) to the table, mirroring the existing
# "2021" column (M) for formatting, then fill in the new figures.
$ws.Range("M2:M6").Copy($ws.Range("N2:N6"))

$ws.Range("N3").Value = 2022
$ws.Range("N4").Value = 6333
$ws.Range("N5").Value = 82675
$ws.Range("N6").Value = 300853

# Match the author's final selection (cell N2) recorded in the sheet view.
[void]$ws.Range("N2").Select()
